$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header "Language Code" in C1, copying the style from the
# existing header cell B1 (bold header style with fill/border).
$ws.Range("C1").Value = "Language Code"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats

# Widen column C to match columns A and B.
$ws.Columns.Item(3).ColumnWidth = 20

# Move the active selection to A2, like in the updated file.
$ws.Range("A2").Select()
